$d = $word.ActiveDocument

function InsertTextAt($doc, $pos, $text) {
    $ins = $doc.Range($pos, $pos)
    $ins.InsertAfter($text)
    return ($pos + $text.Length)
}

# ---------------------------------------------------------------------
# Change 1: fix the doubled "processing" in the code path reference.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    ".../code/processingprocessingcode.R",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ".../code/processing/processingcode.R",
    2)

# ---------------------------------------------------------------------
# Change 2: add a new BodyText paragraph describing the EDA code file,
# right after the "Finally, we examined the most popular sender domain
# names..." paragraph (still inside the same bookmarked section).
# ---------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute(
    "Finally, we examined the most popular sender domain names",
    $true, $false, $false, $false, $false,
    $true, 1, $false, $null, 0)

$srcPara = $rng.Paragraphs(1)
$srcPara.Range.InsertParagraphAfter()

$newPara = $srcPara.Next()
$pos = $newPara.Range.Start

$pos = InsertTextAt $d $pos "A complete set of exploratory visualization steps utilizing"
$pos = InsertTextAt $d $pos " "

$ggStart = $pos
$pos = InsertTextAt $d $pos "ggplot2"
$ggEnd = $pos

$pos = InsertTextAt $d $pos " "
$pos = InsertTextAt $d $pos "can be located in the"
$pos = InsertTextAt $d $pos " "

$edaStart = $pos
$pos = InsertTextAt $d $pos ".../code/eda-code/eda.qmd"
$edaEnd = $pos

$pos = InsertTextAt $d $pos " "
$pos = InsertTextAt $d $pos "file found within the project."

$d.Range($ggStart, $ggEnd).Style = "VerbatimChar"
$d.Range($edaStart, $edaEnd).Style = "VerbatimChar"
